# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") values were regenerated from the re-scraped source data
# (strikeouts per outing, replacing the old "Strike#" derived figure).
# Apply the newly computed K values to the corresponding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    4  = 0
    5  = 2
    6  = 1
    7  = 2
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 0
    16 = 1
    17 = 2
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 2
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 0
    41 = 2
    42 = 2
    43 = 0
    44 = 1
    45 = 0
    46 = 1
    47 = 0
    48 = 0
    49 = 2
    50 = 1
    52 = 2
    53 = 0
    54 = 1
    55 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
